# feat(CWL): mod integrity check
# Adds a "missing mods" warning block (id / text_JP / text columns) as three
# new rows (60-62) at the bottom of the General sheet, and normalizes the
# font style of D59 to match its row's other translated cell (C59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) D59: style-only change (s="10" -> s="3"), value/text stays v=128.
#    Reuse C59's existing format (already style index 3) via PasteSpecial so
#    no new style entry gets minted.
# ---------------------------------------------------------------------------
$ws.Range("C59").Copy()
$ws.Range("D59").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Stamp the formatting for the three new rows (60, 61, 62) up front by
#    cloning row 59's A/C/D formats (A=s2, C=s3, D=s3 after the fix above).
#    This keeps every new cell on an existing style index.
# ---------------------------------------------------------------------------
$ws.Range("A59").Copy()
$ws.Range("A60:A62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C59:D59").Copy()
$ws.Range("C60:D60").PasteSpecial(-4122)
$ws.Range("C61:D61").PasteSpecial(-4122)
$ws.Range("C62:D62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights: row 60 wraps two lines (double height), 61 & 62 are single-line.
$ws.Rows.Item(60).RowHeight = 46.5
$ws.Rows.Item(61).RowHeight = 23.25
$ws.Rows.Item(62).RowHeight = 23.25

# ---------------------------------------------------------------------------
# 3) Values — written in the same order the author entered them so new
#    shared-string indices land at 173..181 in the expected sequence:
#    the three ids first, then each row's JP/ZH text pair.
# ---------------------------------------------------------------------------
$ws.Cells.Item(60,1).Value = "cwl_warn_missing_mods"
$ws.Cells.Item(61,1).Value = "cwl_warn_missing_mods_yes"
$ws.Cells.Item(62,1).Value = "cwl_warn_missing_mods_no"

$ws.Cells.Item(60,3).Value = "現在のセーブから欠落しているMOD：`n{0}"
$ws.Cells.Item(60,4).Value = "当前存档中缺失的模组：`n{0}"

$ws.Cells.Item(61,3).Value = "セーブせずに終了"
$ws.Cells.Item(61,4).Value = "不保存并返回至标题"

$ws.Cells.Item(62,3).Value = "プレイを続ける"
$ws.Cells.Item(62,4).Value = "继续游玩"

# D60 carries mixed formatting: the "{0}" placeholder is rendered in the
# monospace "Cascadia Code" font (matching the other bilingual entries in
# this sheet), the rest uses the sheet's normal font.
# "当前存档中缺失的模组：`n{0}" -> "{0}" starts at (1-based) character 13.
$d60 = $ws.Cells.Item(60,4)
$chars = $d60.Characters(13, 3)
$chars.Font.Name = "Cascadia Code"
$chars.Font.Size = 11

# ---------------------------------------------------------------------------
# 4) Leave the view parked where the author left it: selection on D58.
# ---------------------------------------------------------------------------
$ws.Range("D58").Select() | Out-Null
